$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 158-159, shifting the existing
# rows 158:169 down to 160:171 (matches dimension A1:T169 -> A1:T171).
$ws.Range("A158:A159").EntireRow.Insert()

# --- New row 158 (Primera, $/bandeja 18 kilos) ---
$ws.Cells.Item(158, 1).Value = 7
$ws.Cells.Item(158, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(158, 3).Value = "Ñuble"
$ws.Cells.Item(158, 4).Value = 44753
$ws.Cells.Item(158, 5).Value = 16
$ws.Cells.Item(158, 6).Value = "Fruta"
$ws.Cells.Item(158, 7).Value = 100101
$ws.Cells.Item(158, 8).Value = "Berries"
$ws.Cells.Item(158, 9).Value = 100101007
$ws.Cells.Item(158, 10).Value = "Kiwi"
$ws.Cells.Item(158, 11).Value = "Hayward"
$ws.Cells.Item(158, 12).Value = "Primera"
$ws.Cells.Item(158, 13).Value = 120
$ws.Cells.Item(158, 14).Value = 6500
$ws.Cells.Item(158, 15).Value = 7000
$ws.Cells.Item(158, 16).Value = 6750
$ws.Cells.Item(158, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(158, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(158, 19).Value = 375
$ws.Cells.Item(158, 20).Value = 18

# --- New row 159 (Segunda, $/bandeja 18 kilos) ---
$ws.Cells.Item(159, 1).Value = 7
$ws.Cells.Item(159, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(159, 3).Value = "Ñuble"
$ws.Cells.Item(159, 4).Value = 44753
$ws.Cells.Item(159, 5).Value = 16
$ws.Cells.Item(159, 6).Value = "Fruta"
$ws.Cells.Item(159, 7).Value = 100101
$ws.Cells.Item(159, 8).Value = "Berries"
$ws.Cells.Item(159, 9).Value = 100101007
$ws.Cells.Item(159, 10).Value = "Kiwi"
$ws.Cells.Item(159, 11).Value = "Hayward"
$ws.Cells.Item(159, 12).Value = "Segunda"
$ws.Cells.Item(159, 13).Value = 120
$ws.Cells.Item(159, 14).Value = 5500
$ws.Cells.Item(159, 15).Value = 6000
$ws.Cells.Item(159, 16).Value = 5750
$ws.Cells.Item(159, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(159, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(159, 19).Value = 319
$ws.Cells.Item(159, 20).Value = 18
